# Scheduled "Updated cryptos list" refresh (GitHub Actions bot).
# Rewrites the Price (D) / Volume(1h) (E) columns for each coin row with
# newly scraped figures, and also re-orders the Algorand / WEMIXTOKEN pair
# (rows 36-37) to reflect their new rank.
#
# Values that look numeric but must keep significant trailing/leading
# zeros (e.g. "0.07020", "7.040") are written with a leading apostrophe so
# Excel stores them as literal text instead of silently normalising them
# to a Double and dropping the zero - exactly like typing '0.07020 into
# a cell by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.761.18'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '1.700.46'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").Value = '317.12'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("D7").Value = '0.3932'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '0.4042'
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").Value = '1.505'
$ws.Range("E9").Value = '  -3.04%  '
$ws.Range("D10").Value = '54.04'
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").Value = '0.08895'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '7.249'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '23.43'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '8.023'
$ws.Range("E15").Value = '  +4.90%  '
$ws.Range("D16").Value = '0.00001327'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '1.711.72'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '100.21'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = "'0.07020"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '19.64'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").Value = "'7.040"
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '14.47'
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("D24").Value = '24.761.95'
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = '3.215'
$ws.Range("E25").Value = '  +8.03%  '
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("D28").Value = '161.45'
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("D29").Value = '136.65'
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("D30").Value = '5.174'
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").Value = '7.757'
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("D32").Value = "'0.08720"
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("D33").Value = '1.071'
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("D34").Value = '7.206'
$ws.Range("E34").Value = '  -3.96%  '
$ws.Range("D35").Value = '11.28'
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.957'
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2748'
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("D38").Value = '14.36'
$ws.Range("E38").Value = '  -3.11%  '
$ws.Range("D39").Value = '0.09198'
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("D40").Value = '0.02735'
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("D41").Value = '1.465'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("D42").Value = '0.7695'
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("D43").Value = '16.01'
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("D44").Value = "'0.7190"
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("D45").Value = '2.571'
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("D46").Value = '4.222'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").Value = '140.45'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").Value = '1.308'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").Value = '0.07991'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").Value = '90.51'
$ws.Range("E51").Value = '  +2.50%  '
